# Update the embedded build timestamp throughout the workbook.
# Old build: "January 30 2026 16.19.47 EST"
# New build: "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$oldBuild = "January 30 2026 16.19.47 EST"
$newBuild = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldBuild)"
$newVersion = "mines - January 30 (built on $newBuild)"

$aboutSheet = $wb.Worksheets.Item("About")
$boundariesSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: A2 version banner, A6 recommended citation
$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Mezhegeyugol Coal Mine, Russia, M0820, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: S2:S7 build_version column for each data row
for ($r = 2; $r -le 7; $r++) {
    $cell = $boundariesSheet.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
